$d = $word.ActiveDocument

# The sentence describing sf::Vector3f was missing its closing period.
# "...3 coordinates(x, y and z)" -> "...3 coordinates(x, y and z)."
$d.Content.Find.Execute(
    "coordinates(x, y and z)",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "coordinates(x, y and z).",
    2
)
